$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "1.003"). A leading
# apostrophe is the standard Excel "force text" marker so COM stores the
# literal digits instead of silently parsing them into a rounded Double.
$ws.Range('D2').Value = "28.018.92"
$ws.Range('E2').Value = "  -3.74%  "
$ws.Range('D3').Value = "1.914.03"
$ws.Range('E3').Value = "  -2.99%  "
$ws.Range('D4').Value = "'1.003"
$ws.Range('E4').Value = "  -1.00%  "
$ws.Range('D5').Value = "'329.21"
$ws.Range('E5').Value = "  -0.05%  "
$ws.Range('D6').Value = "'1.003"
$ws.Range('E6').Value = "  -0.96%  "
$ws.Range('D7').Value = "'0.4686"
$ws.Range('E7').Value = "  -5.51%  "
$ws.Range('D8').Value = "'0.4017"
$ws.Range('E8').Value = "  -4.47%  "
$ws.Range('D9').Value = "'53.11"
$ws.Range('E9').Value = "  -2.04%  "
$ws.Range('D10').Value = "'0.08375"
$ws.Range('E10').Value = "  -10.40%  "
$ws.Range('D11').Value = "'1.042"
$ws.Range('E11').Value = "  -5.22%  "
$ws.Range('D12').Value = "'22.05"
$ws.Range('E12').Value = "  -2.94%  "
$ws.Range('B13').Value = "Chainlink"
$ws.Range('C13').Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range('D13').Value = "'7.440"
$ws.Range('E13').Value = "  -5.58%  "
$ws.Range('B14').Value = "WrappedEther"
$ws.Range('C14').Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('D14').Value = "1.865.99"
$ws.Range('E14').Value = "  -9.30%  "
$ws.Range('D15').Value = "'6.048"
$ws.Range('E15').Value = "  -6.28%  "
$ws.Range('D16').Value = "'1.003"
$ws.Range('E16').Value = "  -1.11%  "
$ws.Range('D17').Value = "'89.51"
$ws.Range('E17').Value = "  -2.46%  "
$ws.Range('D18').Value = "'0.00001062"
$ws.Range('E18').Value = "  -4.40%  "
$ws.Range('D19').Value = "'0.06568"
$ws.Range('E19').Value = "  -2.33%  "
$ws.Range('D20').Value = "'17.93"
$ws.Range('E20').Value = "  -6.41%  "
$ws.Range('E21').Value = "  -0.90%  "
$ws.Range('D22').Value = "'5.693"
$ws.Range('E22').Value = "  -4.33%  "
$ws.Range('D23').Value = "28.002.84"
$ws.Range('D24').Value = "'11.32"
$ws.Range('E24').Value = "  -5.36%  "
$ws.Range('D25').Value = "'2.287"
$ws.Range('E25').Value = "  +0.93%  "
$ws.Range('D26').Value = "2.100.64"
$ws.Range('E26').Value = "  -5.32%  "
$ws.Range('D27').Value = "'153.82"
$ws.Range('E27').Value = "  -1.90%  "
$ws.Range('D28').Value = "'19.95"
$ws.Range('E28').Value = "  -3.96%  "
$ws.Range('D29').Value = "'2.125"
$ws.Range('E29').Value = "  -6.11%  "
$ws.Range('D30').Value = "'5.676"
$ws.Range('E30').Value = "  -9.19%  "
$ws.Range('D31').Value = "'122.89"
$ws.Range('E31').Value = "  -3.38%  "
$ws.Range('D32').Value = "'0.9677"
$ws.Range('E32').Value = "  -7.31%  "
$ws.Range('D33').Value = "'0.09563"
$ws.Range('E33').Value = "  -2.68%  "
$ws.Range('D34').Value = "'1.439"
$ws.Range('E34').Value = "  -4.04%  "
$ws.Range('D35').Value = "'3.637"
$ws.Range('E35').Value = "  -2.98%  "
$ws.Range('D36').Value = "'5.511"
$ws.Range('E36').Value = "  -5.14%  "
$ws.Range('D37').Value = "'8.811"
$ws.Range('E37').Value = "  -2.46%  "
$ws.Range('D38').Value = "'0.02295"
$ws.Range('E38').Value = "  -5.05%  "
$ws.Range('D39').Value = "'0.06124"
$ws.Range('E39').Value = "  -4.49%  "
$ws.Range('D40').Value = "'1.213"
$ws.Range('E40').Value = "  -8.53%  "
$ws.Range('D41').Value = "'0.6112"
$ws.Range('E41').Value = "  -5.64%  "
$ws.Range('D42').Value = "'10.98"
$ws.Range('E42').Value = "  -4.43%  "
$ws.Range('E43').Value = "  -0.89%  "
$ws.Range('D44').Value = "'0.1895"
$ws.Range('E44').Value = "  -5.32%  "
$ws.Range('E45').Value = "  -3.71%  "
$ws.Range('D46').Value = "'0.5842"
$ws.Range('D47').Value = "'12.69"
$ws.Range('E47').Value = "  -4.18%  "
$ws.Range('D48').Value = "'2.014"
$ws.Range('E48').Value = "  -7.45%  "
$ws.Range('D49').Value = "'3.465"
$ws.Range('E49').Value = "  -0.42%  "
$ws.Range('D50').Value = "'0.06820"
$ws.Range('E50').Value = "  -2.06%  "
$ws.Range('D51').Value = "'109.33"
$ws.Range('E51').Value = "  -3.41%  "
